$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "long result = 1 + 5 / 1;" "long quantity = 3 + 8 / 1;"
Replace-Text "[6L, 6l]" "[11.0L, 11.0l]"
Replace-Text "double total = 8 * 6;" "double quantity = 3 * 8;"
Replace-Text "[48.0]" "[24.0]"
Replace-Text "int k = 3 / 2;" "int n = 2 / 7;"
Replace-Text "[1]" "[0.2857142857142857]"
Replace-Text "int n = 7 % 4;" "int g = 3 % 6;"
Replace-Text "[3]" "[3.0]"
Replace-Text "double m = 7 / 4 * (61 - 6);" "double n = 7.0 / 5 * (53 - 2);"
Replace-Text "[55.0]" "[71.39999999999999]"
